$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.795758356961096
$ws.Range("D2").Value = 4.501983508384387
$ws.Range("E2").Value = 12.53507889104725
$ws.Range("F2").Value = 23.33724731771041
$ws.Range("G2").Value = 3.619230112743357
$ws.Range("K2").Value = 9.912723146394059
$ws.Range("M2").Value = 13.55153538810813
$ws.Range("N2").Value = 18.43473513515739
$ws.Range("O2").Value = 20.77647784399818

$ws.Range("B3").Value = 7.723603566338228
$ws.Range("D3").Value = 4.458564282370052
$ws.Range("E3").Value = 12.38627156865967
$ws.Range("F3").Value = 23.29445716053176
$ws.Range("G3").Value = 3.621195632881195
$ws.Range("K3").Value = 9.496848281443539
$ws.Range("M3").Value = 13.31738852825255
$ws.Range("N3").Value = 18.49753026609397
$ws.Range("O3").Value = 20.80250768729889

$ws.Range("B4").Value = 7.680808492176975
$ws.Range("D4").Value = 4.431274144845615
$ws.Range("E4").Value = 12.29855333332908
$ws.Range("F4").Value = 23.27508232687383
$ws.Range("G4").Value = 3.622466752269673
$ws.Range("K4").Value = 9.229959984299446
$ws.Range("M4").Value = 13.17478351248557
$ws.Range("N4").Value = 18.53791511055742
$ws.Range("O4").Value = 20.82375350563835

$ws.Range("B5").Value = 7.663767917852403
$ws.Range("D5").Value = 4.419999469347108
$ws.Range("E5").Value = 12.26376859896418
$ws.Range("F5").Value = 23.26892688972621
$ws.Range("G5").Value = 3.62300095971539
$ws.Range("K5").Value = 9.118388169401504
$ws.Range("M5").Value = 13.11704217195276
$ws.Range("N5").Value = 18.55483340782999
$ws.Range("O5").Value = 20.83373201904863

$ws.Range("B6").Value = 7.660962983182275
$ws.Range("D6").Value = 4.418118176399584
$ws.Range("E6").Value = 12.25805189023635
$ws.Range("F6").Value = 23.26800999152356
$ws.Range("G6").Value = 3.623090645381522
$ws.Range("K6").Value = 9.099694736323261
$ws.Range("M6").Value = 13.10747901113649
$ws.Range("N6").Value = 18.557670570993
$ws.Range("O6").Value = 20.83546861469725

$ws.Range("B7").Value = 7.680577037225672
$ws.Range("D7").Value = 4.431122706613849
$ws.Range("E7").Value = 12.29808026710126
$ws.Range("F7").Value = 23.27499226188082
$ws.Range("G7").Value = 3.622473891051497
$ws.Range("K7").Value = 9.228466548936028
$ws.Range("M7").Value = 13.17400318427787
$ws.Range("N7").Value = 18.53814140764146
$ws.Range("O7").Value = 20.82388273622487

$ws.Range("B8").Value = 7.77057826646847
$ws.Range("D8").Value = 4.487146304202549
$ws.Range("E8").Value = 12.48303926363864
$ws.Range("F8").Value = 23.32106521437436
$ws.Range("G8").Value = 3.619894513522195
$ws.Range("K8").Value = 9.771773956822036
$ws.Range("M8").Value = 13.47060773790357
$ws.Range("N8").Value = 18.45600825226729
$ws.Range("O8").Value = 20.78435896273734

$ws.Range("B9").Value = 7.95814785245336
$ws.Range("D9").Value = 4.591798527143538
$ws.Range("E9").Value = 12.87269275108407
$ws.Range("F9").Value = 23.4658344615291
$ws.Range("G9").Value = 3.615344051486098
$ws.Range("K9").Value = 10.74235347512884
$ws.Range("M9").Value = 14.05802068658871
$ws.Range("N9").Value = 18.30938873098125
$ws.Range("O9").Value = 20.74872061237275

$ws.Range("B10").Value = 8.101479118919164
$ws.Range("D10").Value = 4.665242299752196
$ws.Range("E10").Value = 13.17258628549963
$ws.Range("F10").Value = 23.60482762707103
$ws.Range("G10").Value = 3.612307025150626
$ws.Range("K10").Value = 11.39399127268876
$ws.Range("M10").Value = 14.48858365965502
$ws.Range("N10").Value = 18.21038155990865
$ws.Range("O10").Value = 20.74817075458794

$ws.Range("B11").Value = 8.16761535891194
$ws.Range("D11").Value = 4.697851208964387
$ws.Range("E11").Value = 13.31135238677304
$ws.Range("F11").Value = 23.6749900705512
$ws.Range("G11").Value = 3.610991188610317
$ws.Range("K11").Value = 11.67651773339391
$ws.Range("M11").Value = 14.68331981083919
$ws.Range("N11").Value = 18.16721348015057
$ws.Range("O11").Value = 20.75349763505927

$ws.Range("B12").Value = 8.192771096145673
$ws.Range("D12").Value = 4.710079684650294
$ws.Range("E12").Value = 13.36418409639056
$ws.Range("F12").Value = 23.70253993639062
$ws.Range("G12").Value = 3.610502312991732
$ws.Range("K12").Value = 11.78146522863109
$ws.Range("M12").Value = 14.75682606482599
$ws.Range("N12").Value = 18.15113443556923
$ws.Range("O12").Value = 20.75631650533369

$ws.Range("B13").Value = 8.187348769013406
$ws.Range("D13").Value = 4.70745146438863
$ws.Range("E13").Value = 13.35279405135633
$ws.Range("F13").Value = 23.69656323146821
$ws.Range("G13").Value = 3.610607183651323
$ws.Range("K13").Value = 11.75895415552898
$ws.Range("M13").Value = 14.74100681714878
$ws.Range("N13").Value = 18.15458545604627
$ws.Range("O13").Value = 20.7556737669898

$ws.Range("B14").Value = 8.169682835545068
$ws.Range("D14").Value = 4.698859680200016
$ws.Range("E14").Value = 13.31569349380141
$ws.Range("F14").Value = 23.67723705535138
$ws.Range("G14").Value = 3.610950780338217
$ws.Range("K14").Value = 11.68519289627184
$ws.Range("M14").Value = 14.68937235695859
$ws.Range("N14").Value = 18.16588528804792
$ws.Range("O14").Value = 20.75371348196876

$ws.Range("B15").Value = 8.158875766063133
$ws.Range("D15").Value = 4.69358122652991
$ws.Range("E15").Value = 13.29300374937963
$ws.Range("F15").Value = 23.66552644277265
$ws.Range("G15").Value = 3.611162466387129
$ws.Range("K15").Value = 11.63974538443769
$ws.Range("M15").Value = 14.65771187996195
$ws.Range("N15").Value = 18.17284159628364
$ws.Range("O15").Value = 20.75261713579518

$ws.Range("B16").Value = 8.097173794359952
$ws.Range("D16").Value = 4.66309471378517
$ws.Range("E16").Value = 13.1635602137216
$ws.Range("F16").Value = 23.600380526569
$ws.Range("G16").Value = 3.61239433657784
$ws.Range("K16").Value = 11.37524414342158
$ws.Range("M16").Value = 14.47582809608918
$ws.Range("N16").Value = 18.21324023354732
$ws.Range("O16").Value = 20.74793484627826

$ws.Range("B17").Value = 8.059544547599041
$ws.Range("D17").Value = 4.644183730705555
$ws.Range("E17").Value = 13.08471320335129
$ws.Range("F17").Value = 23.56218054997474
$ws.Range("G17").Value = 3.613166848014634
$ws.Range("K17").Value = 11.2093891347043
$ws.Range("M17").Value = 14.363905608093
$ws.Range("N17").Value = 18.23850175022244
$ws.Range("O17").Value = 20.74649083952683

$ws.Range("B18").Value = 8.037990499217038
$ws.Range("D18").Value = 4.633231504784063
$ws.Range("E18").Value = 13.03958619545557
$ws.Range("F18").Value = 23.54086265358126
$ws.Range("G18").Value = 3.613617365046457
$ws.Range("K18").Value = 11.11268663026921
$ws.Range("M18").Value = 14.29942869144196
$ws.Range("N18").Value = 18.25320766664291
$ws.Range("O18").Value = 20.74618527837068

$ws.Range("B19").Value = 8.030708705879301
$ws.Range("D19").Value = 4.62951050892039
$ws.Range("E19").Value = 13.02434700309589
$ws.Range("F19").Value = 23.53375752210704
$ws.Range("G19").Value = 3.61377096677656
$ws.Range("K19").Value = 11.07972164107708
$ws.Range("M19").Value = 14.27758268598993
$ws.Range("N19").Value = 18.25821713111106
$ws.Range("O19").Value = 20.74617198607004

$ws.Range("B20").Value = 8.063541169440212
$ws.Range("D20").Value = 4.646204646575821
$ws.Range("E20").Value = 13.0930838333959
$ws.Range("F20").Value = 23.56617945692691
$ws.Range("G20").Value = 3.613083972672654
$ws.Range("K20").Value = 11.22718026889713
$ws.Range("M20").Value = 14.37583105553766
$ws.Range("N20").Value = 18.23579439577966
$ws.Range("O20").Value = 20.7465902214816

$ws.Range("B21").Value = 8.174868908172977
$ws.Range("D21").Value = 4.701386583279818
$ws.Range("E21").Value = 13.32658353976165
$ws.Range("F21").Value = 23.68288713877796
$ws.Range("G21").Value = 3.610849602841613
$ws.Range("K21").Value = 11.70691398858341
$ws.Range("M21").Value = 14.70454562844496
$ws.Range("N21").Value = 18.1625589961284
$ws.Range("O21").Value = 20.75426751288204

$ws.Range("B22").Value = 8.248266654198282
$ws.Range("D22").Value = 4.736750773826136
$ws.Range("E22").Value = 13.48082017301282
$ws.Range("F22").Value = 23.76487091451761
$ws.Range("G22").Value = 3.60944409974167
$ws.Range("K22").Value = 12.00854827687281
$ws.Range("M22").Value = 14.91797162589663
$ws.Range("N22").Value = 18.11625563050745
$ws.Range("O22").Value = 20.76395736461268

$ws.Range("B23").Value = 8.209042017116612
$ws.Range("D23").Value = 4.717941818030478
$ws.Range("E23").Value = 13.39836933544101
$ws.Range("F23").Value = 23.72059808771222
$ws.Range("G23").Value = 3.61018924584494
$ws.Range("K23").Value = 11.84866050115439
$ws.Range("M23").Value = 14.80421447525617
$ws.Range("N23").Value = 18.14082625781317
$ws.Range("O23").Value = 20.75835845631118

$ws.Range("B24").Value = 8.061734047310685
$ws.Range("D24").Value = 4.645291239488554
$ws.Range("E24").Value = 13.08929883462101
$ws.Range("F24").Value = 23.56436954452618
$ws.Range("G24").Value = 3.613121420706342
$ws.Range("K24").Value = 11.21914109986721
$ws.Range("M24").Value = 14.37043996370634
$ws.Range("N24").Value = 18.23701782140357
$ws.Range("O24").Value = 20.74654365670009

$ws.Range("B25").Value = 7.906344746016217
$ws.Range("D25").Value = 4.564073165145445
$ws.Range("E25").Value = 12.76468169466614
$ws.Range("F25").Value = 23.42089332170319
$ws.Range("G25").Value = 3.616521063184012
$ws.Range("K25").Value = 10.49035258770289
$ws.Range("M25").Value = 13.89896664708248
$ws.Range("N25").Value = 18.34751637589879
$ws.Range("O25").Value = 20.75386572189163

